$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking.com "cryptos" sheet: refresh Price (D) / Volume(1h) (E)
# columns for the latest snapshot. Both columns are plain text cells
# (values like "44.138.70" or "0.830" are not valid numbers, and the
# percentages carry significant leading/trailing padding spaces), so
# Set-TextValue forces text formatting before the write and clears the
# format override again afterwards to avoid leaving stray number formats
# behind on cells that had none originally.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "44.138.70"
$ws.Range("E2").Value = "  +0.49%  "
Set-TextValue $ws.Range("D3") "2.246.76"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "306.63"
$ws.Range("E5").Value = "  -2.71%  "
Set-TextValue $ws.Range("D6") "96.32"
$ws.Range("E6").Value = "  -3.41%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  +0.18%  "
Set-TextValue $ws.Range("D9") "0.529"
$ws.Range("E9").Value = "  -1.36%  "
Set-TextValue $ws.Range("D10") "35.13"
$ws.Range("E10").Value = "  -3.57%  "
Set-TextValue $ws.Range("D11") "0.0813"
$ws.Range("E11").Value = "  -1.10%  "
Set-TextValue $ws.Range("D12") "7.27"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("E13").Value = "  -0.23%  "
Set-TextValue $ws.Range("D14") "2.588.65"
$ws.Range("E14").Value = "  +0.45%  "
Set-TextValue $ws.Range("D15") "2.334.30"
$ws.Range("E15").Value = "  +4.46%  "
Set-TextValue $ws.Range("D16") "0.830"
$ws.Range("E16").Value = "  -1.32%  "
Set-TextValue $ws.Range("D17") "13.63"
$ws.Range("E17").Value = "  -2.88%  "
Set-TextValue $ws.Range("D18") "44.017.20"
$ws.Range("E18").Value = "  +0.46%  "
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  -5.01%  "
Set-TextValue $ws.Range("D21") "6.35"
$ws.Range("E21").Value = "  -0.21%  "
Set-TextValue $ws.Range("D22") "65.55"
$ws.Range("E22").Value = "  +0.91%  "
Set-TextValue $ws.Range("D23") "237.31"
$ws.Range("E23").Value = "  +1.42%  "
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("E25").Value = "  -2.62%  "
Set-TextValue $ws.Range("D27") "38.82"
$ws.Range("E27").Value = "  +4.97%  "
Set-TextValue $ws.Range("D28") "9.94"
$ws.Range("E28").Value = "  -3.09%  "
$ws.Range("E29").Value = "  +0.96%  "
Set-TextValue $ws.Range("D30") "6.10"
$ws.Range("E30").Value = "  +0.73%  "
Set-TextValue $ws.Range("D31") "20.06"
$ws.Range("E31").Value = "  +0.45%  "
Set-TextValue $ws.Range("D32") "151.69"
$ws.Range("E32").Value = "  -4.23%  "
Set-TextValue $ws.Range("D33") "0.0807"
$ws.Range("E33").Value = "  -3.49%  "
Set-TextValue $ws.Range("D34") "3.31"
$ws.Range("E34").Value = "  +4.72%  "
$ws.Range("E35").Value = "  -2.75%  "
$ws.Range("E36").Value = "  +0.09%  "
Set-TextValue $ws.Range("D37") "0.121"
$ws.Range("E37").Value = "  +2.64%  "
Set-TextValue $ws.Range("D38") "1.79"
$ws.Range("E38").Value = "  -4.85%  "
Set-TextValue $ws.Range("D39") "15.23"
$ws.Range("E39").Value = "  -4.28%  "
Set-TextValue $ws.Range("D40") "3.44"
$ws.Range("E40").Value = "  -5.05%  "
Set-TextValue $ws.Range("D41") "3.86"
$ws.Range("E41").Value = "  -4.10%  "
$ws.Range("E42").Value = "  -3.15%  "
$ws.Range("E43").Value = "  +0.13%  "
Set-TextValue $ws.Range("D44") "1.739.52"
$ws.Range("E44").Value = "  +0.63%  "
Set-TextValue $ws.Range("D45") "84.86"
$ws.Range("E45").Value = "  +5.08%  "
Set-TextValue $ws.Range("D46") "0.190"
$ws.Range("E46").Value = "  -2.49%  "
Set-TextValue $ws.Range("D47") "100.58"
Set-TextValue $ws.Range("D48") "4.98"
$ws.Range("E48").Value = "  -2.58%  "
Set-TextValue $ws.Range("D49") "70.04"
$ws.Range("E49").Value = "  -4.75%  "
Set-TextValue $ws.Range("D50") "8.12"
$ws.Range("E50").Value = "  +0.30%  "
Set-TextValue $ws.Range("D51") "54.42"
$ws.Range("E51").Value = "  -4.35%  "
